$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.072.32"
$ws.Range("E2").Value = "  -2.90%  "
$ws.Range("D3").Value = "1.866.34"
$ws.Range("E3").Value = "  -2.21%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'306.49"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").Value = "'0.5097"
$ws.Range("E7").Value = "  +1.76%  "
$ws.Range("D8").Value = "'0.3743"
$ws.Range("E8").Value = "  -1.64%  "
$ws.Range("D9").Value = "'0.07132"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").Value = "'0.8873"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D11").Value = "'20.62"
$ws.Range("E11").Value = "  -2.91%  "
$ws.Range("D12").Value = "1.876.52"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "'0.07518"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "'5.309"
$ws.Range("E14").Value = "  -3.07%  "
$ws.Range("D15").Value = "'89.01"
$ws.Range("E15").Value = "  -3.92%  "
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'0.000008475"
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("E18").Value = "  -3.58%  "
$ws.Range("D19").Value = "'0.9992"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("D20").Value = "27.127.40"
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "'5.044"
$ws.Range("E21").Value = "  -2.52%  "
$ws.Range("D22").Value = "2.104.96"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("D23").Value = "'10.54"
$ws.Range("E23").Value = "  -2.96%  "
$ws.Range("D24").Value = "'6.473"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("D25").Value = "'149.56"
$ws.Range("E25").Value = "  -2.38%  "
$ws.Range("D26").Value = "'1.847"
$ws.Range("E26").Value = "  +0.63%  "
$ws.Range("D27").Value = "'17.93"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "'2.099"
$ws.Range("E28").Value = "  -4.79%  "
$ws.Range("D29").Value = "'112.75"
$ws.Range("E29").Value = "  -2.14%  "
$ws.Range("D30").Value = "'4.726"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("D31").Value = "'4.683"
$ws.Range("E31").Value = "  -2.89%  "
$ws.Range("D32").Value = "'0.09019"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("D33").Value = "'0.05125"
$ws.Range("E33").Value = "  -2.62%  "
$ws.Range("D34").Value = "'3.084"
$ws.Range("E34").Value = "  -3.78%  "
$ws.Range("E35").Value = "  -5.92%  "
$ws.Range("D36").Value = "'0.7337"
$ws.Range("E36").Value = "  -5.71%  "
$ws.Range("D37").Value = "'0.02048"
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").Value = "'2.513"
$ws.Range("E38").Value = "  -4.39%  "
$ws.Range("D39").Value = "'3.060"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("D40").Value = "'1.078"
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").Value = "'0.5314"
$ws.Range("E41").Value = "  -4.12%  "
$ws.Range("D42").Value = "'6.560"
$ws.Range("D43").Value = "'116.59"
$ws.Range("E43").Value = "  +2.10%  "
$ws.Range("D44").Value = "'8.324"
$ws.Range("E44").Value = "  -2.29%  "
$ws.Range("D45").Value = "'0.1468"
$ws.Range("E45").Value = "  -3.30%  "
$ws.Range("D46").Value = "'0.9987"
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("D47").Value = "'0.4622"
$ws.Range("E47").Value = "  -4.01%  "
$ws.Range("D48").Value = "'10.05"
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("D49").Value = "'1.567"
$ws.Range("E49").Value = "  -4.11%  "
$ws.Range("D50").Value = "'64.40"
$ws.Range("E50").Value = "  -4.40%  "
$ws.Range("D51").Value = "'36.44"
$ws.Range("E51").Value = "  -1.77%  "
